# Change document/table header attributes to lowerCamelCase across sheets.
$wb = $excel.ActiveWorkbook

$wsLeaves = $wb.Worksheets.Item("!!Leaves")
$wsLeaves.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$wsLeaves.Range("A2").Value = "!!ObjTables type='Data' id='Leaf'"

$wsNormal = $wb.Worksheets.Item("!!Normal records")
$wsNormal.Range("A1").Value = "!!ObjTables type='Data' id='NormalRecord'"

$wsTransposed = $wb.Worksheets.Item("!!Transposed")
$wsTransposed.Range("A1").Value = "!!ObjTables type='Data' id='Transposed'"
